# Work Profile and new tenant support
#
# Appends new interview-history rows to the AMSIN, BETA and AMS sheets of
# the mass-interview-history tracker, and fixes the run-time timestamp /
# formatting of the previously-added AMS "live173" row (row 44).
#
# Columns (all four history sheets): A=Run Date (text, e.g. "2023-03-10"),
# B=Run Time (date/time serial), C=Sprint Name (text), D=Total Cases,
# E=Pass Cases, F=Fail Cases, G=Time Taken.

$wb = $excel.ActiveWorkbook

$xlPasteValues  = -4163
$xlPasteFormats = -4122

# Column A holds the run date as literal text that happens to look like a
# date (e.g. "2023-03-10"). A plain `.Value = "2023-03-10"` assignment
# gets auto-converted by Excel into a real date. Routing the literal
# through a tiny `="..."` formula and then collapsing it back to a static
# value with Copy/PasteSpecial(values) keeps it a plain text cell without
# leaving a stray quote-prefix behind.
function Set-TextCell {
    param($cell, [string]$text)
    $cell.ClearContents()
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

# Plain numeric cell (counts, durations, ...).
function Set-NumberCell {
    param($cell, [double]$number)
    $cell.ClearContents()
    $cell.Value = $number
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

# Column B holds the actual run-time timestamp (date+time serial number)
# and needs the sheet's "YYYY-MM-DD HH:MM:SS" number format. Copy that
# format from a known-good cell in the same column after writing the
# value (assigning .Value first clears any pasted formatting).
function Set_RunTimeCell {
    param($ws, $cell, [double]$number, [string]$formatSourceA1)
    $cell.ClearContents()
    $cell.Value = $number
    $ws.Range($formatSourceA1).Copy()
    $cell.PasteSpecial($xlPasteFormats)
}

function Add-HistoryRow {
    param(
        $ws,
        [int]$row,
        [string]$runDate,
        [double]$runTime,
        [string]$sprintName,
        [double]$total,
        [double]$pass,
        [double]$fail,
        [double]$timeTaken,
        [string]$formatSourceA1
    )

    Set-TextCell   ($ws.Cells.Item($row, 1)) $runDate
    Set_RunTimeCell $ws ($ws.Cells.Item($row, 2)) $runTime $formatSourceA1
    Set-TextCell   ($ws.Cells.Item($row, 3)) $sprintName
    Set-NumberCell ($ws.Cells.Item($row, 4)) $total
    Set-NumberCell ($ws.Cells.Item($row, 5)) $pass
    Set-NumberCell ($ws.Cells.Item($row, 6)) $fail
    Set-NumberCell ($ws.Cells.Item($row, 7)) $timeTaken
}

# ---------------------------------------------------------------------
# AMSIN sheet: append rows 61-65
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Add-HistoryRow $wsAmsin 61 "2023-03-10" 44995.7960968287   "174ffiinnalrun" 155 154 1 4.1  "B60"
Add-HistoryRow $wsAmsin 62 "2023-03-13" 44998.47009619213  "174finalrun"    155 155 0 2.81 "B61"
Add-HistoryRow $wsAmsin 63 "2023-03-28" 45013.53435076389  "175prerun"      155 155 0 3.46 "B62"
Add-HistoryRow $wsAmsin 64 "2023-03-31" 45016.4500043287   "175fnlrun"      155 155 0 2.98 "B63"
Add-HistoryRow $wsAmsin 65 "2023-04-12" 45028.60502657553  "176fstrtail"    155 154 1 3.39 "B64"

# ---------------------------------------------------------------------
# BETA sheet: append rows 33-34
# ---------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

Add-HistoryRow $wsBeta 33 "2023-03-13" 44998.53603851852 "174beta" 155 149 6 4.17 "B32"
Add-HistoryRow $wsBeta 34 "2023-03-31" 45016.54302363426 "175beta" 155 154 1 3.3  "B33"

# ---------------------------------------------------------------------
# AMS sheet: fix existing row 44 (live173) and append rows 45-47
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Row 44 keeps its original text/numbers; only the run-time timestamp and
# the cell styling (it was missing the table's usual style) change.
Set-TextCell    ($wsAms.Cells.Item(44, 1)) "2023-02-20"
Set_RunTimeCell $wsAms ($wsAms.Cells.Item(44, 2)) 44977.84034873843 "B43"
Set-TextCell    ($wsAms.Cells.Item(44, 3)) "live173"
Set-NumberCell  ($wsAms.Cells.Item(44, 4)) 155
Set-NumberCell  ($wsAms.Cells.Item(44, 5)) 147
Set-NumberCell  ($wsAms.Cells.Item(44, 6)) 8
Set-NumberCell  ($wsAms.Cells.Item(44, 7)) 4.27

Add-HistoryRow $wsAms 45 "2023-03-01" 44986.69101804398 "173angularvrs" 155 153 2 3.19 "B44"
Add-HistoryRow $wsAms 46 "2023-03-13" 44998.83911592593 "174live"       155 151 4 3.8  "B45"
Add-HistoryRow $wsAms 47 "2023-03-31" 45016.80636225695 "175live"       155 152 3 3.44 "B46"
